# Actualización desde MV -datos-
# Appends three new daily rows (04-10-2021, 05-10-2021, 06-10-2021) to the
# "Spot posiciones netas y suscripciones" sheet, right after the existing
# last data row (01-10-2021, row 190).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Date = "04-10-2021"; B = -10111; C = 2595; D = 673;  E = 742;  F = 1180 },
    @{ Date = "05-10-2021"; B = -9946;  C = 3743; D = 1094; E = 974;  F = 1675 },
    @{ Date = "06-10-2021"; B = -10060; C = 3321; D = 973;  E = 1075; F = 1273 }
)

$startRow = 191
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    # Column A holds the date label as TEXT (shared string), matching the
    # existing "Serie" column. Force text interpretation by temporarily
    # switching the cell to a text number format, then restore the default
    # ("Normal") style so the cell ends up stored just like its neighbours.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data.Date
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 3).Value = $data.C
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
}
